$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; temporarily unprotect so the cells can be edited.
$ws.Unprotect()

# Update the "as of" date in the confidentiality disclaimer text (A10).
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."
# Re-fit the row height (the embedded line break can otherwise stamp an
# explicit row height even though the wrapped text still spans two lines).
$ws.Rows(10).AutoFit()

# Update Weight (D) and Percent Change (E) values for the holdings rows.
$ws.Range("D2").Value = 0.4880583390011705
$ws.Range("E2").Value = 0.003100775193798366

$ws.Range("D3").Value = 0.3313783702888761
$ws.Range("E3").Value = 0.007909383849233453

$ws.Range("D4").Value = 0.09692717145575572
$ws.Range("E4").Value = 0.003595182455509605

$ws.Range("D5").Value = 0.05422455587076958
$ws.Range("E5").Value = 0.002639430801010034

$ws.Range("D6").Value = 0.0294115633834281
$ws.Range("E6").Value = 0.021580450650587

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.005260665542018428

# Restore sheet protection to its prior (protected) state.
$ws.Protect()
